$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 0.02354566666666667
$ws.Range("H2").Value = 0.070637
$ws.Range("I2").Value = 0.002815555392485919
$ws.Range("J2").Value = 0.002815555392485918
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 7.423863
$ws.Range("N2").Value = 22.271589
$ws.Range("O2").Value = 0.1690720838224332
$ws.Range("P2").Value = 0.1690720838224332
$ws.Range("Q2").Value = 0.174799803577
$ws.Range("R2").Value = 1.573198232193
$ws.Range("S2").Value = 0.0004760318173250832
$ws.Range("T2").Value = 0.0004760318173250831

# Row 3
$ws.Range("G3").Value = 0.02354566666666667
$ws.Range("H3").Value = 0.070637
$ws.Range("I3").Value = 0.002815555392485919
$ws.Range("J3").Value = 0.002815555392485918
$ws.Range("O3").Value = 0.4908369772207905
$ws.Range("P3").Value = 0.4908369772207905
$ws.Range("Q3").Value = 0.507465249536
$ws.Range("R3").Value = 4.567187245824
$ws.Range("S3").Value = 0.001381978698045485
$ws.Range("T3").Value = 0.001381978698045485

# Row 4
$ws.Range("G4").Value = 0.02354566666666667
$ws.Range("H4").Value = 0.070637
$ws.Range("I4").Value = 0.002815555392485919
$ws.Range("J4").Value = 0.002815555392485918
$ws.Range("O4").Value = 0.3400909389567762
$ws.Range("P4").Value = 0.3400909389567762
$ws.Range("Q4").Value = 0.3516123299834445
$ws.Range("R4").Value = 3.164510969851
$ws.Range("S4").Value = 0.0009575448771153507
$ws.Range("T4").Value = 0.0009575448771153506

# Row 5
$ws.Range("I5").Value = 0.9868456480383168
$ws.Range("J5").Value = 0.9868456480383166
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 7.423863
$ws.Range("N5").Value = 22.271589
$ws.Range("O5").Value = 0.1690720838224332
$ws.Range("P5").Value = 0.1690720838224332
$ws.Range("Q5").Value = 61.266926553205
$ws.Range("R5").Value = 551.402338978845
$ws.Range("S5").Value = 0.1668480501249378
$ws.Range("T5").Value = 0.1668480501249377

# Row 6
$ws.Range("I6").Value = 0.9868456480383168
$ws.Range("J6").Value = 0.9868456480383166
$ws.Range("O6").Value = 0.4908369772207905
$ws.Range("P6").Value = 0.4908369772207905
$ws.Range("S6").Value = 0.4843803348666195
$ws.Range("T6").Value = 0.4843803348666195

# Row 7
$ws.Range("I7").Value = 0.9868456480383168
$ws.Range("J7").Value = 0.9868456480383166
$ws.Range("O7").Value = 0.3400909389567762
$ws.Range("P7").Value = 0.3400909389567762
$ws.Range("S7").Value = 0.3356172630467595
$ws.Range("T7").Value = 0.3356172630467594

# Row 8
$ws.Range("I8").Value = 0.0103387965691973
$ws.Range("J8").Value = 0.0103387965691973
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 7.423863
$ws.Range("N8").Value = 22.271589
$ws.Range("O8").Value = 0.1690720838224332
$ws.Range("P8").Value = 0.1690720838224332
$ws.Range("Q8").Value = 0.6418696696009999
$ws.Range("R8").Value = 5.776827026408999
$ws.Range("S8").Value = 0.001748001880170412
$ws.Range("T8").Value = 0.001748001880170412

# Row 9
$ws.Range("I9").Value = 0.0103387965691973
$ws.Range("J9").Value = 0.0103387965691973
$ws.Range("O9").Value = 0.4908369772207905
$ws.Range("P9").Value = 0.4908369772207905
$ws.Range("S9").Value = 0.005074663656125484
$ws.Range("T9").Value = 0.005074663656125484

# Row 10
$ws.Range("I10").Value = 0.0103387965691973
$ws.Range("J10").Value = 0.0103387965691973
$ws.Range("O10").Value = 0.3400909389567762
$ws.Range("P10").Value = 0.3400909389567762
$ws.Range("S10").Value = 0.003516131032901407
$ws.Range("T10").Value = 0.003516131032901407
